# Append " (Changed main)" to the end of the first paragraph
# ("This is a Microsoft word document.") as three additional runs:
#   " (", "Changed main", ")"
#
# A plain Range.InsertAfter() would coalesce each new piece of text into
# the run that already sits at the insertion point (same, empty,
# formatting => same run). To get three *separate* <w:r> elements - as
# the target OOXML expects - insert the text and then immediately wrap it
# in a temporary bookmark and delete the bookmark again. Adding the
# bookmark forces Word to split the run at that boundary; removing the
# bookmark only strips the <w:bookmarkStart/</w:bookmarkEnd> markers and
# leaves the run split in place (with no leftover formatting residue).

$d = $word.ActiveDocument

$para = $d.Paragraphs(1).Range
$paraEnd = $para.End - 1   # position right before the paragraph mark

function Insert-AsSeparateRun($start, $text, $bookmarkName) {
    $r = $d.Range($start, $start)
    $r.InsertAfter($text)
    $d.Bookmarks.Add($bookmarkName, $r)
    $d.Bookmarks($bookmarkName).Delete()
    return $start + $text.Length
}

$pos = $paraEnd
$pos = Insert-AsSeparateRun $pos " (" "__tmp_split_1"
$pos = Insert-AsSeparateRun $pos "Changed main" "__tmp_split_2"
$pos = Insert-AsSeparateRun $pos ")" "__tmp_split_3"
